$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend header row (row 1) with two new columns P1, Q1 continuing the sequence
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15

# Apply the same formatting as the other header cells (bold font, thin border,
# centered horizontally, top-aligned vertically) to the new header cells, by
# copying the format from the adjacent existing header cell O1
$ws.Range("O1").Copy()
$ws.Range("P1:Q1").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

# Update existing I, K, M, O columns (swap 1<->2) for rows 2-25
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 9).Value = 2   # I column
    $ws.Cells.Item($r, 11).Value = 1  # K column
    $ws.Cells.Item($r, 13).Value = 2  # M column
    $ws.Cells.Item($r, 15).Value = 1  # O column
    $ws.Cells.Item($r, 16).Value = 2  # P column (new)
    $ws.Cells.Item($r, 17).Value = 2  # Q column (new)
}
